$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "67.226.93"
$ws.Range("E2").Value = "  -2.88%  "
Set-TextValue "D3" "3.777.19"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue "D5" "592.70"
$ws.Range("E5").Value = "  -1.27%  "
Set-TextValue "D6" "165.44"
$ws.Range("E6").Value = "  -3.19%  "
Set-TextValue "D7" "3.782.50"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  -4.04%  "
Set-TextValue "D14" "35.76"
$ws.Range("E14").Value = "  -3.01%  "
Set-TextValue "D15" "4.410.57"
$ws.Range("E15").Value = "  -0.71%  "
Set-TextValue "D16" "3.764.09"
$ws.Range("E16").Value = "  -0.90%  "
Set-TextValue "D17" "67.278.06"
$ws.Range("E17").Value = "  -2.84%  "
Set-TextValue "D18" "17.92"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  -1.82%  "
Set-TextValue "D21" "10.20"
$ws.Range("E21").Value = "  -7.89%  "
Set-TextValue "D22" "456.52"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("E24").Value = "  +0.37%  "
Set-TextValue "D25" "83.21"
$ws.Range("E25").Value = "  -1.96%  "
Set-TextValue "D26" "11.81"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("E27").Value = "  -5.47%  "
$ws.Range("E28").Value = "  +0.01%  "
Set-TextValue "D29" "9.91"
$ws.Range("E29").Value = "  -3.58%  "
$ws.Range("E30").Value = "  -2.08%  "
Set-TextValue "D31" "29.73"
$ws.Range("E31").Value = "  -2.04%  "
Set-TextValue "D32" "2.19"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("E33").Value = "  -4.06%  "
Set-TextValue "D34" "9.15"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("E35").Value = "  +0.13%  "
Set-TextValue "D36" "3.727.69"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("E39").Value = "  -7.21%  "
$ws.Range("E40").Value = "  -1.78%  "
Set-TextValue "D41" "5.71"
$ws.Range("E41").Value = "  -2.91%  "
Set-TextValue "D42" "1.00"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  -0.01%  "
Set-TextValue "D44" "43.65"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("E45").Value = "  -4.19%  "
Set-TextValue "D46" "46.90"
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("E47").Value = "  -3.71%  "
Set-TextValue "D48" "147.45"
$ws.Range("E48").Value = "  +1.68%  "
Set-TextValue "D49" "392.08"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("E50").Value = "  -8.07%  "
Set-TextValue "D51" "2.744.45"
$ws.Range("E51").Value = "  +1.78%  "
